# "made simulink model more clear and updated some parameters"
#
# The motor-selection table at the top of the sheet (rows 1-5) gets a new
# "rotor damping (N*m*s)" column, the "stall torque" / "no load speed"
# columns swap order, a new "Maximum current" value is filled in for the
# lab motor, and the three commercial-motor rows are reordered. Further
# down, a handful of blank styled filler cells shift from column F to a
# newly-used column G (to keep the blank "spacer" column between the data
# tables and the adjoining helper tables consistent), and the last little
# results table (rows 47-52) moves its "Speed(rpm)" column from F to G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-Str($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = $text
}
function Set-Num($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}
function Clear-Cell($row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

# ---------------------------------------------------------------------
# Row 1: header row — add "rotor damping (N*m*s)" (new col F), shift the
# old "Maximum current (A)" header out to col G, and swap the stall
# torque / no-load-speed headers (D <-> E).
# ---------------------------------------------------------------------
Set-Str 1 4 "stall torque (Nm)"
Set-Str 1 5 "no load speed (rpm)"
Set-Str 1 6 "rotor damping (N*m*s)"
Set-Str 1 7 "Maximum current (A)"
$ws.Cells.Item(1, 6).Font.Bold = $true
$ws.Cells.Item(1, 7).Font.Bold = $true

# ---------------------------------------------------------------------
# Row 2: GM8724S009 (Lab1) — swap D/E (stall torque, no load speed),
# add rotor damping + max current values.
# ---------------------------------------------------------------------
Set-Num 2 4 0.3
Set-Num 2 5 720
Set-Num 2 6 0.0000014
$ws.Cells.Item(2, 6).NumberFormat = "0.00E+00"
Set-Num 2 7 2.77

# ---------------------------------------------------------------------
# Rows 3-5: the three commercial motors get reordered
#   old row3 = E30-150-48, old row4 = E30-400-48, old row5 = QBL5704
#   new row3 = QBL5704,    new row4 = E30-400-48, new row5 = E30-150-48
# plus rotor damping values are added for QBL5704 and E30-400-48, and
# the E30-400-48 / E30-150-48 "armature inductance" values switch from
# the placeholder text "0 (unlisted)" (E30-400-48 only) to real numbers.
# ---------------------------------------------------------------------

# New row 3 (was row 5: QBL5704-116-04-042)
Set-Str 3 1 "QBL5704-116-04-042"
Set-Num 3 2 36
Set-Num 3 3 0.001
$ws.Cells.Item(3, 3).NumberFormat = "0.00E+00"
Set-Num 3 4 1.3
Set-Num 3 5 5500
Set-Num 3 6 0.000001
$ws.Cells.Item(3, 6).NumberFormat = "0.00E+00"
Set-Num 3 7 11

# New row 4 (was row 4: E30-400-48) - stays same motor, gains real
# armature inductance + rotor damping values (was "0 (unlisted)").
Set-Str 4 1 "E30-400-48"
Set-Num 4 2 45
Set-Num 4 3 0.0025
$ws.Cells.Item(4, 3).NumberFormat = "0.00E+00"
Set-Num 4 4 13.13
Set-Num 4 5 5800
Set-Num 4 6 0.000002
$ws.Cells.Item(4, 6).NumberFormat = "0.00E+00"
Set-Num 4 7 280

# New row 5 (was row 3: E30-150-48) - still "0 (unlisted)" armature
# inductance, no rotor damping value given.
Set-Str 5 1 "E30-150-48 "
Set-Num 5 2 45
Set-Str 5 3 "0 (unlisted)"
Set-Num 5 4 5.3
Set-Num 5 5 5300
Clear-Cell 5 6
Set-Num 5 7 70

# ---------------------------------------------------------------------
# Blank styled "spacer" cells: a handful of header/divider rows pick up
# an extra blank bold cell in column F (matching the rest of the row),
# while stray blank styled cells that used to sit in column G move out
# (no longer needed at G once F is in use) for rows 31/35/38/39.
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 6).Font.Bold = $true
$ws.Cells.Item(20, 6).Font.Bold = $true
$ws.Cells.Item(32, 6).Font.Bold = $true
$ws.Cells.Item(44, 6).Font.Bold = $true

Clear-Cell 31 7
Clear-Cell 35 7
Clear-Cell 38 7
Clear-Cell 39 7

# Row 47: the stray bold/yellow-fill cell at G47 becomes a plain bold
# centered cell (matching the merged C47:G47 header band style).
$ws.Cells.Item(47, 7).Font.Bold = $true
$ws.Cells.Item(47, 7).HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Row 48 + results table (rows 49-52): the "Speed(rpm)" column and its
# data move from column F to column G.
# ---------------------------------------------------------------------
Set-Str 48 7 "Speed(rpm)"
$ws.Cells.Item(48, 7).Font.Bold = $true
Clear-Cell 48 6
$ws.Cells.Item(48, 6).Font.Bold = $true

Set-Num 49 7 256.7445
Set-Num 50 7 436.4015
Set-Num 51 7 580.8943
Set-Num 52 7 703.1914
Clear-Cell 49 6
Clear-Cell 50 6
Clear-Cell 51 6
Clear-Cell 52 6

# ---------------------------------------------------------------------
# Merged cells: the result-summary header merge and the secondary-table
# label merge both grow by one column (to keep pace with the new G
# column usage).
# ---------------------------------------------------------------------
$ws.Range("C47:F47").UnMerge()
$ws.Range("C47:G47").Merge()

$ws.Range("G31:J31").UnMerge()
$ws.Range("H31:J31").Merge()

# ---------------------------------------------------------------------
# Column widths: new column F inherits column E's width; old column F's
# width now belongs to column G.
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth
$ws.Columns("G").ColumnWidth = 23.5703125

# Selection, as last left by the editing user.
$ws.Range("A3").Select()
